# Atualização de bases das ligas, do dia: 02-05-2024 às 20:28
# This script swaps the data of row 60 <-> row 61, and performs a cyclic
# rotation of row 139 -> 141 -> 142 -> 139 (i.e. row 141 gets what was in
# row 139, row 142 gets what was in row 141, and row 139 gets what was in
# row 142), leaving the id column (A) and the Div/Date columns (C/D)
# untouched since they are identical across the rows involved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B..AB hold the data that needs to move between rows.
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB")

function Get-RowValues($row) {
    $vals = @{}
    foreach ($col in $cols) {
        $vals[$col] = $ws.Range("$col$row").Value()
    }
    return $vals
}

function Set-RowValues($row, $vals) {
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value = $vals[$col]
    }
}

# --- Swap row 60 and row 61 ---
$row60 = Get-RowValues 60
$row61 = Get-RowValues 61

Set-RowValues 60 $row61
Set-RowValues 61 $row60

# --- Cyclic rotation among rows 139, 141, 142 ---
# After edit: row139 <= old row142 ; row141 <= old row139 ; row142 <= old row141
$row139 = Get-RowValues 139
$row141 = Get-RowValues 141
$row142 = Get-RowValues 142

Set-RowValues 139 $row142
Set-RowValues 141 $row139
Set-RowValues 142 $row141
